# Ran scraper to update the data: insert "height" and "weight" columns
# between the existing "fumbles" and "fantasy points" columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at E:F, shifting the existing "fantasy
# points" column (old E) to G. This also carries the header style (s="1")
# from the old E1 cell onto the two new header cells, matching the diff.
$ws.Columns("E:F").Insert()

# New header labels for the inserted columns.
$ws.Cells.Item(1, 5).Value = "height"
$ws.Cells.Item(1, 6).Value = "weight"

# Fill the new height/weight columns for every data row (2-17). The
# scraped values are constant across all rows in this sheet.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 255
}
